# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, new DAMSLTag (col I), new DialogAct (col J)
$updates = @(
    ,@(7, 'b', 'Acknowledge (Backchannel)')
    ,@(8, 'sd', 'Statement-non-opinion')
    ,@(17, 'b', 'Acknowledge (Backchannel)')
    ,@(22, 'sd', 'Statement-non-opinion')
    ,@(24, 'sv', 'Statement-opinion')
    ,@(39, 'sd', 'Statement-non-opinion')
    ,@(41, 'sd', 'Statement-non-opinion')
    ,@(45, 'sd', 'Statement-non-opinion')
    ,@(46, 'sd', 'Statement-non-opinion')
    ,@(51, 'sd', 'Statement-non-opinion')
    ,@(63, '%', 'Uninterpretable')
    ,@(64, '%', 'Uninterpretable')
    ,@(67, 'aa', 'Agree/Accept')
    ,@(85, 'b', 'Acknowledge (Backchannel)')
    ,@(87, 'sd', 'Statement-non-opinion')
    ,@(88, 'sd', 'Statement-non-opinion')
    ,@(93, 'sd', 'Statement-non-opinion')
    ,@(98, 'sd', 'Statement-non-opinion')
    ,@(117, 'ba', 'Appreciation')
    ,@(130, 'b', 'Acknowledge (Backchannel)')
    ,@(138, 'sd', 'Statement-non-opinion')
    ,@(139, 'sd', 'Statement-non-opinion')
    ,@(141, 'b', 'Acknowledge (Backchannel)')
    ,@(147, 'sd', 'Statement-non-opinion')
    ,@(150, 'sv', 'Statement-opinion')
    ,@(154, 'sv', 'Statement-opinion')
    ,@(155, 'aa', 'Agree/Accept')
    ,@(178, 'qy', 'Yes-No-Question')
    ,@(195, 'b', 'Acknowledge (Backchannel)')
    ,@(197, 'sd', 'Statement-non-opinion')
    ,@(202, 'aa', 'Agree/Accept')
    ,@(204, 'aa', 'Agree/Accept')
    ,@(231, 'sv', 'Statement-opinion')
    ,@(245, 'b', 'Acknowledge (Backchannel)')
    ,@(248, '%', 'Uninterpretable')
    ,@(252, '%', 'Uninterpretable')
    ,@(260, 'b', 'Acknowledge (Backchannel)')
    ,@(264, 'b', 'Acknowledge (Backchannel)')
    ,@(271, 'sv', 'Statement-opinion')
    ,@(272, 'sd', 'Statement-non-opinion')
    ,@(273, 'sd', 'Statement-non-opinion')
    ,@(278, 'ba', 'Appreciation')
    ,@(287, 'sd', 'Statement-non-opinion')
    ,@(292, 'sd', 'Statement-non-opinion')
    ,@(317, 'b', 'Acknowledge (Backchannel)')
    ,@(320, 'sd', 'Statement-non-opinion')
    ,@(321, 'b', 'Acknowledge (Backchannel)')
    ,@(332, 'b', 'Acknowledge (Backchannel)')
    ,@(334, '%', 'Uninterpretable')
    ,@(340, 'sd', 'Statement-non-opinion')
    ,@(341, 'sd', 'Statement-non-opinion')
    ,@(342, 'aa', 'Agree/Accept')
    ,@(350, 'b', 'Acknowledge (Backchannel)')
    ,@(359, 'aa', 'Agree/Accept')
    ,@(365, 'sv', 'Statement-opinion')
    ,@(401, 'sd', 'Statement-non-opinion')
    ,@(413, 'sd', 'Statement-non-opinion')
    ,@(414, 'sd', 'Statement-non-opinion')
    ,@(415, 'sd', 'Statement-non-opinion')
    ,@(423, 'sv', 'Statement-opinion')
    ,@(425, 'aa', 'Agree/Accept')
    ,@(435, 'aa', 'Agree/Accept')
    ,@(447, 'sd', 'Statement-non-opinion')
    ,@(459, 'sd', 'Statement-non-opinion')
    ,@(466, 'b', 'Acknowledge (Backchannel)')
    ,@(476, 'b', 'Acknowledge (Backchannel)')
    ,@(483, 'sv', 'Statement-opinion')
    ,@(491, 'sd', 'Statement-non-opinion')
    ,@(495, 'sd', 'Statement-non-opinion')
    ,@(498, 'b', 'Acknowledge (Backchannel)')
    ,@(502, '%', 'Uninterpretable')
    ,@(507, 'ba', 'Appreciation')
    ,@(516, 'sd', 'Statement-non-opinion')
    ,@(525, 'ba', 'Appreciation')
    ,@(536, 'aa', 'Agree/Accept')
    ,@(547, 'b', 'Acknowledge (Backchannel)')
    ,@(552, 'sv', 'Statement-opinion')
    ,@(555, 'sd', 'Statement-non-opinion')
    ,@(561, 'sd', 'Statement-non-opinion')
    ,@(569, 'sd', 'Statement-non-opinion')
    ,@(570, 'sd', 'Statement-non-opinion')
    ,@(577, 'sd', 'Statement-non-opinion')
    ,@(579, 'sv', 'Statement-opinion')
    ,@(582, 'b', 'Acknowledge (Backchannel)')
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 9).Value = $u[1]
    $ws.Cells.Item($r, 10).Value = $u[2]
}
